$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fill in missing / updated "Лаба №4" (column G) scores for several students.
$ws.Range("G3").Value = 0
$ws.Range("G10").Value = 5
$ws.Range("G11").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("G24").Value = 5
$ws.Range("G25").Value = -1
$ws.Range("G26").Value = 5

# A small, empty formatted block was added below the table (rows 32-35).
$ws.Range("B32:L32").Value = ""
$ws.Range("B33:K35").Value = ""

# Restore the on-screen selection to the newly added block.
$ws.Range("B32:K35").Select()
